$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rows where column N (Main Number Partner Asset) should become the
# "intra receiver main" scenario value.
$mainRows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)
foreach ($r in $mainRows) {
    $ws.Cells.Item($r, 14).Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_MAIN"
}

# Rows where column N currently holds the numeric placeholder 126 and
# should become the "intra receiver sub main" scenario value.
$subMainRows = @(8, 13, 18, 23, 28)
foreach ($r in $subMainRows) {
    $ws.Cells.Item($r, 14).Value = "MPA_TRANSFER_SCENARIO_INTRA_REC_SUB_MAIN"
}
